# Restore cell C10 on the "Rules" sheet back to its saved value of 1
# (previously 18) as recorded by revision #a775ad755117bb96b664b97d63b242eba4039ddf.TEST

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
